# Weekly update: insert a new price-report row at row 7 (right after the
# most recent existing rows), pushing the previous rows 7-76 down to 8-77.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a new row at position 7.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(7, 1).Value  = 7
$ws.Cells.Item(7, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value  = "Ñuble"
$ws.Cells.Item(7, 4).Value  = 45061
$ws.Cells.Item(7, 5).Value  = 16
$ws.Cells.Item(7, 6).Value  = 100112001
$ws.Cells.Item(7, 7).Value  = "Berenjena"
$ws.Cells.Item(7, 8).Value  = "Sin especificar"
$ws.Cells.Item(7, 9).Value  = "Primera"
$ws.Cells.Item(7, 10).Value = 60
$ws.Cells.Item(7, 11).Value = 8000
$ws.Cells.Item(7, 12).Value = 9000
$ws.Cells.Item(7, 13).Value = 8500
$ws.Cells.Item(7, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 142
$ws.Cells.Item(7, 17).Value = 60
$ws.Cells.Item(7, 18).Value = "Hortaliza"
